$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header and value for the new "visibility_level" column (M)
$ws.Range("M1").Value = "visibility_level"
$ws.Range("M2").Value = "PRO"

# Set the width of the new column M to match the target width (18.83203125),
# closest achievable value given COM ColumnWidth quantization.
$ws.Columns.Item(13).ColumnWidth = 17.9986979166667
